# "add Crisis and Credit Allocation"
# The DAG regression-output table is refreshed for the new sample/
# specification: the "A" variable (and its lag) is replaced by "FFR"
# (Fed Funds Rate), shifting the column order to C / FFR / LF, and all
# of the reported coefficients are updated to the new estimates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (Source / C / FFR / LF) ---
$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "FFR"
$ws.Range("D1").Value = "LF"

# --- Row labels (left column) ---
$ws.Range("A2").Value = "C Lag"
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("A4").Value = "LF Lag"

# --- Coefficient values ---
# These look numeric but must stay plain text (to match the asterisk-
# significance-flagged entries elsewhere in the table), so force the
# cells to Text format before writing them.
$coeffs = $ws.Range("B2:D4")
$coeffs.NumberFormat = "@"

$ws.Range("B2").Value = "-0.46***"
$ws.Range("C2").Value = "3.79"
$ws.Range("D2").Value = "-6.09"

$ws.Range("B3").Value = "-0.01"
$ws.Range("C3").Value = "1.6***"
$ws.Range("D3").Value = "0.5***"

$ws.Range("B4").Value = "0.04*"
$ws.Range("C4").Value = "3.53*"
$ws.Range("D4").Value = "0.54*"

# Drop back to the default style so we don't leave stray explicit
# formatting on cells that were unstyled before this edit.
$coeffs.Style = "Normal"
